$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.188.59"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "2.317.49"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D5").Value = "'302.29"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "'97.73"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").Value = "'35.46"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "'19.46"
$ws.Range("E11").Value = "  +7.07%  "
$ws.Range("D12").Value = "'0.0794"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "2.684.13"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "2.322.32"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "'0.790"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "43.131.84"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'12.60"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").Value = "'6.05"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "'67.93"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'236.91"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "'2.24"
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").Value = "'2.44"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'24.98"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "  +7.90%  "
$ws.Range("D29").Value = "'164.30"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Value = "'9.11"
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("D31").Value = "'32.95"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'17.91"
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("D34").Value = "'5.00"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").Value = "'4.48"
$ws.Range("E35").Value = "  -6.94%  "
$ws.Range("D36").Value = "'0.0701"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D37").Value = "'2.36"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").Value = "'1.77"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").Value = "1.984.56"
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").Value = "'10.69"
$ws.Range("E43").Value = "  +6.50%  "
$ws.Range("D44").Value = "'19.02"
$ws.Range("E44").Value = "  +7.35%  "
$ws.Range("D45").Value = "'0.0280"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").Value = "2.550.42"
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("D50").Value = "'53.83"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "'72.54"
$ws.Range("E51").Value = "  +0.60%  "

# Reset style on cells that needed a text quote-prefix so no stray number format / style is left behind
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
